$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row data: Row, NumDoc (C), NombreTrabajador (D), PeriodoMora (E), ValorMora (F), SalarioBasico (G)
$rows = @(
    @(16, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1805", 31249, 1684984),
    @(17, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1806", 31249, 1684984),
    @(18, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1807", 31249, 1684984),
    @(19, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1808", 31249, 1684984),
    @(20, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1809", 31249, 1684984),
    @(21, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1810", 31249, 1684984),
    @(22, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1811", 31249, 1684984),
    @(23, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1812", 31249, 1684984),
    @(24, "1047365908", "SHIRLEY PAOLA DIAZ MARRUGO", "1812", 37697, 1684984),
    @(25, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1901", 31249, 1684984),
    @(26, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1902", 31249, 1684984),
    @(27, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1903", 31249, 1684984),
    @(28, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1904", 31249, 1684984),
    @(29, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1905", 31249, 1684984),
    @(30, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1906", 31249, 1684984),
    @(31, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1907", 31249, 1684984),
    @(32, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1908", 31249, 1684984),
    @(33, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1909", 31249, 1684984),
    @(34, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1910", 31249, 1684984),
    @(35, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1911", 31249, 1684984),
    @(36, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "1912", 31249, 1684984),
    @(37, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2001", 31249, 1684984),
    @(38, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2002", 31249, 1684984),
    @(39, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2003", 31249, 1684984),
    @(40, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2004", 31249, 1684984),
    @(41, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2005", 31249, 1684984),
    @(42, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2006", 31249, 1684984),
    @(43, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2007", 31249, 1684984),
    @(44, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2008", 31249, 1684984),
    @(45, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2009", 31249, 1684984),
    @(46, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2010", 31249, 1684984),
    @(47, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2011", 31249, 1684984),
    @(48, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2012", 31249, 1684984),
    @(49, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2101", 31249, 1684984),
    @(50, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2102", 31249, 1684984),
    @(51, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2103", 31249, 1684984),
    @(52, "1047427900", "JOSE JAVIER ROMERO MANJARRES", "2104", 31249, 1684984)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
}

Write-Host "Done updating rows 16-52"
